$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1071
$ws.Range("I18").Value = 979.4
$ws.Range("J18").Value = 1300
$ws.Range("K18").Value = 979.4
$ws.Range("L18").Value = 1300
$ws.Range("M18").Value = -695.4
$ws.Range("N18").Value = -1868

$ws.Range("H33").Value = 878
$ws.Range("I33").Value = 1004.7
$ws.Range("K33").Value = 1004.7
$ws.Range("M33").Value = -775.7

$ws.Range("H41").Value = 2337
$ws.Range("J41").Value = 2174.5
$ws.Range("L41").Value = 2174.5
$ws.Range("N41").Value = -3054.5

$ws.Range("H80").Value = 2693.5625
$ws.Range("I80").Value = 194.16667
$ws.Range("J80").Value = 4193.2
$ws.Range("K80").Value = 582.50001
$ws.Range("L80").Value = 12579.6
$ws.Range("M80").Value = 415.49999
$ws.Range("N80").Value = -14575.6

$ws.Range("H83").Value = 2693.5625
$ws.Range("I83").Value = 194.16667
$ws.Range("J83").Value = 4193.2
$ws.Range("K83").Value = 1747.50003
$ws.Range("L83").Value = 37738.8
$ws.Range("M83").Value = 3244.49997
$ws.Range("N83").Value = -47722.8

$ws.Range("H86").Value = 3999.75
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 4333
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 4333
$ws.Range("M86").Value = -1877
$ws.Range("N86").Value = -6579

$ws.Range("H89").Value = 3999.75
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 4333
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 21665
$ws.Range("M89").Value = -9384
$ws.Range("N89").Value = -32897

$ws.Range("H98").Value = 2493.4614
$ws.Range("I98").Value = 954.25
$ws.Range("J98").Value = 3177.5557
$ws.Range("K98").Value = 954.25
$ws.Range("L98").Value = 3177.5557
$ws.Range("M98").Value = 543.75
$ws.Range("N98").Value = -6173.5557

$ws.Range("H122").Value = 2493.4614
$ws.Range("I122").Value = 954.25
$ws.Range("J122").Value = 3177.5557
$ws.Range("K122").Value = 2862.75
$ws.Range("L122").Value = 9532.667099999999
$ws.Range("M122").Value = -412.75
$ws.Range("N122").Value = -14432.6671

$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -60120

$ws.Range("H135").Value = 220.6
$ws.Range("J135").Value = 518
$ws.Range("L135").Value = 4662
$ws.Range("N135").Value = -9732

$ws.Range("H137").Value = 2039
$ws.Range("I137").Value = 1906.7693
$ws.Range("J137").Value = 2284.5715
$ws.Range("K137").Value = 5720.3079
$ws.Range("L137").Value = 6853.7145
$ws.Range("M137").Value = -3170.3079
$ws.Range("N137").Value = -11953.7145

$ws.Range("H141").Value = 1506.9524
$ws.Range("I141").Value = 1506.9524
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4520.857199999999
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 659.1428000000005
$ws.Range("N141").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1847.36
$ws.Range("I74").Value = 1838.4348
$ws.Range("J74").Value = 1950
$ws.Range("K74").Value = 1838.4348
$ws.Range("L74").Value = 1950
$ws.Range("M74").Value = -964.4348
$ws.Range("N74").Value = -3698

$ws.Range("H77").Value = 1847.36
$ws.Range("I77").Value = 1838.4348
$ws.Range("J77").Value = 1950
$ws.Range("K77").Value = 9192.173999999999
$ws.Range("L77").Value = 9750
$ws.Range("M77").Value = -4824.173999999999
$ws.Range("N77").Value = -18486

$ws.Range("H80").Value = 80000
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""

$ws.Range("H83").Value = 80000
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""

$ws.Range("H122").Value = 2130.6
$ws.Range("I122").Value = 845.375
$ws.Range("K122").Value = 2536.125
$ws.Range("M122").Value = -86.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 90.52941
$ws.Range("I7").Value = 45.666668
$ws.Range("K7").Value = 45.666668
$ws.Range("M7").Value = 67.333332

$ws.Range("H31").Value = 2060
$ws.Range("I31").Value = 1571.875
$ws.Range("K31").Value = 1571.875
$ws.Range("M31").Value = -1276.875

$ws.Range("H34").Value = 2060
$ws.Range("I34").Value = 1571.875
$ws.Range("K34").Value = 1571.875
$ws.Range("M34").Value = -1369.875

$ws.Range("H58").Value = 1769.2609
$ws.Range("I58").Value = 1869.579
$ws.Range("K58").Value = 1869.579
$ws.Range("M58").Value = -1666.579

$ws.Range("H132").Value = 3433.111
$ws.Range("I132").Value = 3899.6667
$ws.Range("K132").Value = 11699.0001
$ws.Range("M132").Value = -9169.000100000001

$ws.Range("H136").Value = 1769.2609
$ws.Range("I136").Value = 1869.579
$ws.Range("K136").Value = 5608.737
$ws.Range("M136").Value = -3058.737

$ws.Range("H141").Value = 207141.28
$ws.Range("J141").Value = 234333
$ws.Range("L141").Value = 234333
$ws.Range("N141").Value = -244693

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6876.1816
$ws.Range("I70").Value = 6571
$ws.Range("K70").Value = 6571
$ws.Range("M70").Value = -6301

$ws.Range("H73").Value = 6876.1816
$ws.Range("I73").Value = 6571
$ws.Range("K73").Value = 6571
$ws.Range("M73").Value = -5635

$ws.Range("H103").Value = 30000
$ws.Range("J103").Value = 30000
$ws.Range("L103").Value = 30000
$ws.Range("N103").Value = -32344

$ws.Range("H122").Value = 2489.1
$ws.Range("I122").Value = 1799.1428
$ws.Range("J122").Value = 4099
$ws.Range("K122").Value = 5397.428400000001
$ws.Range("L122").Value = 12297
$ws.Range("M122").Value = -2947.428400000001
$ws.Range("N122").Value = -17197

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = ""

$ws.Range("H132").Value = 1672.7
$ws.Range("I132").Value = 1672.7
$ws.Range("K132").Value = 5018.1
$ws.Range("M132").Value = -2488.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8532.666999999999
$ws.Range("I7").Value = 6000
$ws.Range("K7").Value = 6000
$ws.Range("M7").Value = -5888

$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = ""

$ws.Range("H40").Value = 2991.6667
$ws.Range("I40").Value = 2321
$ws.Range("J40").Value = 4333
$ws.Range("K40").Value = 2321
$ws.Range("L40").Value = 4333
$ws.Range("M40").Value = -2185
$ws.Range("N40").Value = -4605

$ws.Range("H46").Value = 2475.6428
$ws.Range("J46").Value = 2914.4285
$ws.Range("L46").Value = 2914.4285
$ws.Range("N46").Value = -3290.4285

$ws.Range("H64").Value = 1600
$ws.Range("I64").Value = 1600
$ws.Range("K64").Value = 1600
$ws.Range("M64").Value = -1375

$ws.Range("H67").Value = 1600
$ws.Range("I67").Value = 1600
$ws.Range("K67").Value = 1600
$ws.Range("M67").Value = -820

$ws.Range("H122").Value = 6446.6665
$ws.Range("I122").Value = 8822.727999999999
$ws.Range("K122").Value = 26468.184
$ws.Range("M122").Value = -24018.184

$ws.Range("H126").Value = 8532.666999999999
$ws.Range("I126").Value = 6000
$ws.Range("K126").Value = 18000
$ws.Range("M126").Value = -15530

$ws.Range("H132").Value = 2286.7646
$ws.Range("I132").Value = 1991.2142
$ws.Range("J132").Value = 3666
$ws.Range("K132").Value = 5973.642599999999
$ws.Range("L132").Value = 10998
$ws.Range("M132").Value = -3443.642599999999
$ws.Range("N132").Value = -16058

$ws.Range("H136").Value = 3597.5386
$ws.Range("I136").Value = 2788.3
$ws.Range("K136").Value = 8364.900000000001
$ws.Range("M136").Value = -5814.900000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 8062.5
$ws.Range("I63").Value = 4000
$ws.Range("J63").Value = 12125
$ws.Range("K63").Value = 4000
$ws.Range("L63").Value = 12125
$ws.Range("M63").Value = -3376
$ws.Range("N63").Value = -13373

$ws.Range("H66").Value = 8062.5
$ws.Range("I66").Value = 4000
$ws.Range("J66").Value = 12125
$ws.Range("K66").Value = 12000
$ws.Range("L66").Value = 36375
$ws.Range("M66").Value = -8880
$ws.Range("N66").Value = -42615

$ws.Range("H126").Value = 4958
$ws.Range("I126").Value = 5223.75
$ws.Range("K126").Value = 15671.25
$ws.Range("M126").Value = -13201.25

$ws.Range("H132").Value = 3871.3076
$ws.Range("I132").Value = 4283.7
$ws.Range("K132").Value = 12851.1
$ws.Range("M132").Value = -10321.1

$ws.Range("H136").Value = 2928.9534
$ws.Range("I136").Value = 2965.225
$ws.Range("J136").Value = 2445.3333
$ws.Range("K136").Value = 8895.674999999999
$ws.Range("L136").Value = 7335.999899999999
$ws.Range("M136").Value = -6345.674999999999
$ws.Range("N136").Value = -12435.9999
